$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027695581954037
$ws.Range("D2").Value = 1.030082177297474
$ws.Range("E2").Value = 1.040823199875672
$ws.Range("F2").Value = 1.047573600140512
$ws.Range("I2").Value = 1.029505176095935
$ws.Range("J2").Value = 1.032852024241253
$ws.Range("K2").Value = 1.032894220226891
$ws.Range("L2").Value = 1.043604456557567
$ws.Range("M2").Value = 1.050335859101388
$ws.Range("N2").Value = 1.01494183392994

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028650610735627
$ws.Range("D3").Value = 1.030746753482891
$ws.Range("E3").Value = 1.04179940450206
$ws.Range("F3").Value = 1.048693428204322
$ws.Range("I3").Value = 1.029615614210677
$ws.Range("J3").Value = 1.033447318807099
$ws.Range("K3").Value = 1.033367463180291
$ws.Range("L3").Value = 1.044390727685705
$ws.Range("M3").Value = 1.051266760406308
$ws.Range("N3").Value = 1.015141772147094

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029268972686818
$ws.Range("D4").Value = 1.031176825105245
$ws.Range("E4").Value = 1.042431892356023
$ws.Range("F4").Value = 1.049419154658682
$ws.Range("I4").Value = 1.02968570376308
$ws.Range("J4").Value = 1.033832312183432
$ws.Range("K4").Value = 1.033673048623762
$ws.Range("L4").Value = 1.04489969954149
$ws.Range("M4").Value = 1.051869650735621
$ws.Range("N4").Value = 1.015271003186126

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029529025784786
$ws.Range("D5").Value = 1.031357637021648
$ws.Range("E5").Value = 1.042697985186934
$ws.Range("F5").Value = 1.049724518095813
$ws.Range("I5").Value = 1.029714840631138
$ws.Range("J5").Value = 1.033994114277601
$ws.Range("K5").Value = 1.033801363975613
$ws.Range("L5").Value = 1.045113718821274
$ws.Range("M5").Value = 1.052123233117728
$ws.Range("N5").Value = 1.015325297514849

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029572695320482
$ws.Range("D6").Value = 1.031387996675242
$ws.Range("E6").Value = 1.042742674748123
$ws.Range("F6").Value = 1.049775805656082
$ws.Range("I6").Value = 1.029719713545466
$ws.Range("J6").Value = 1.034021278658656
$ws.Range("K6").Value = 1.033822899703392
$ws.Range("L6").Value = 1.045149656362579
$ws.Range("M6").Value = 1.052165818157196
$ws.Range("N6").Value = 1.015334411745021

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029272447161601
$ws.Range("D7").Value = 1.031179241086028
$ws.Range("E7").Value = 1.042435447136653
$ws.Range("F7").Value = 1.04942323388763
$ws.Range("I7").Value = 1.029686094383828
$ws.Range("J7").Value = 1.033834474384896
$ws.Range("K7").Value = 1.033674763781174
$ws.Range("L7").Value = 1.044902559091909
$ws.Range("M7").Value = 1.051873038616636
$ws.Range("N7").Value = 1.015271728805173

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02801825630115
$ws.Range("D8").Value = 1.030306763544401
$ws.Range("E8").Value = 1.04115294282484
$ws.Range("F8").Value = 1.04795181832762
$ws.Range("I8").Value = 1.029542782814695
$ws.Range("J8").Value = 1.03305324819886
$ws.Range("K8").Value = 1.033054285716609
$ws.Range("L8").Value = 1.043870138280315
$ws.Range("M8").Value = 1.050650350325105
$ws.Range("N8").Value = 1.015009433270415

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025811263867779
$ws.Range("D9").Value = 1.028769759212278
$ws.Range("E9").Value = 1.038899314006773
$ws.Range("F9").Value = 1.045367629419752
$ws.Range("I9").Value = 1.029279768647451
$ws.Range("J9").Value = 1.031675114075158
$ws.Range("K9").Value = 1.031956104775161
$ws.Range("L9").Value = 1.042052462766449
$ws.Range("M9").Value = 1.048499948575118
$ws.Range("N9").Value = 1.0145461588094

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024342027679922
$ws.Range("D10").Value = 1.027745443541265
$ws.Range("E10").Value = 1.037401190001204
$ws.Range("F10").Value = 1.043650685244589
$ws.Range("I10").Value = 1.029097407923331
$ws.Range("J10").Value = 1.030755384296514
$ws.Range("K10").Value = 1.031220796263795
$ws.Range("L10").Value = 1.040841789747584
$ws.Range("M10").Value = 1.047069169693483
$ws.Range("N10").Value = 1.014236604229287

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023706336294948
$ws.Range("D11").Value = 1.027302003388417
$ws.Range("E11").Value = 1.036753515524648
$ws.Range("F11").Value = 1.042908625301583
$ws.Range("I11").Value = 1.029016785109059
$ws.Range("J11").Value = 1.030356909172383
$ws.Range("K11").Value = 1.030901654232805
$ws.Range("L11").Value = 1.040317827918445
$ws.Range("M11").Value = 1.046450304974459
$ws.Range("N11").Value = 1.01410240030326

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023470287532005
$ws.Range("D12").Value = 1.02713730554232
$ws.Range("E12").Value = 1.036513094929485
$ws.Range("F12").Value = 1.042633200394811
$ws.Range("I12").Value = 1.028986589347116
$ws.Range("J12").Value = 1.030208864484796
$ws.Range("K12").Value = 1.030782999037159
$ws.Range("L12").Value = 1.040123246081136
$ws.Range("M12").Value = 1.046220532678753
$ws.Range("N12").Value = 1.014052526531361

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023520917375348
$ws.Range("D13").Value = 1.027172633064661
$ws.Range("E13").Value = 1.036564658956461
$ws.Range("F13").Value = 1.042692270495483
$ws.Range("I13").Value = 1.028993077694722
$ws.Range("J13").Value = 1.03024062206811
$ws.Range("K13").Value = 1.030808456017795
$ws.Range("L13").Value = 1.040164982682745
$ws.Range("M13").Value = 1.046269815007211
$ws.Range("N13").Value = 1.014063225726814

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023686822890739
$ws.Range("D14").Value = 1.027288389088195
$ws.Range("E14").Value = 1.036733639125788
$ws.Range("F14").Value = 1.04288585432165
$ws.Range("I14").Value = 1.029014294192001
$ws.Range("J14").Value = 1.030344672427384
$ws.Range("K14").Value = 1.030891848429837
$ws.Range("L14").Value = 1.040301742884264
$ws.Range("M14").Value = 1.046431309848059
$ws.Range("N14").Value = 1.014098278218962

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023789052752201
$ws.Range("D15").Value = 1.02735971229549
$ws.Range("E15").Value = 1.036837773888606
$ws.Range("F15").Value = 1.043005155445581
$ws.Range("I15").Value = 1.029027333413439
$ws.Range("J15").Value = 1.030408776897653
$ws.Range("K15").Value = 1.030943214479289
$ws.Range("L15").Value = 1.040386010801205
$ws.Range("M15").Value = 1.046530825634378
$ws.Range("N15").Value = 1.014119871984064

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024384226955169
$ws.Range("D16").Value = 1.027774875318589
$ws.Range("E16").Value = 1.037444195649838
$ws.Range("F16").Value = 1.043699962621072
$ws.Range("I16").Value = 1.029102723670346
$ws.Range("J16").Value = 1.030781825073216
$ws.Range("K16").Value = 1.03124196096836
$ws.Range("L16").Value = 1.040876569067071
$ws.Range("M16").Value = 1.047110255920652
$ws.Range("N16").Value = 1.014245507455826

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024757697622718
$ws.Range("D17").Value = 1.028035322513187
$ws.Range("E17").Value = 1.037824862500639
$ws.Range("F17").Value = 1.044136168915417
$ws.Range("I17").Value = 1.029149570030892
$ws.Range("J17").Value = 1.031015768293822
$ws.Range("K17").Value = 1.031429156916086
$ws.Range("L17").Value = 1.041184355448884
$ws.Range("M17").Value = 1.047473897602463
$ws.Range("N17").Value = 1.014324271366968

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02497558471552
$ws.Range("D18").Value = 1.028187246065703
$ws.Range("E18").Value = 1.038046997692996
$ws.Range("F18").Value = 1.044390734566794
$ws.Range("I18").Value = 1.029176734580522
$ws.Range("J18").Value = 1.031152201559228
$ws.Range("K18").Value = 1.031538272861098
$ws.Range("L18").Value = 1.041363907807295
$ws.Range("M18").Value = 1.047686068434913
$ws.Range("N18").Value = 1.014370197108086

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025049886703034
$ws.Range("D19").Value = 1.028239049578131
$ws.Range("E19").Value = 1.038122756762012
$ws.Range("F19").Value = 1.044477557567921
$ws.Range("I19").Value = 1.029185969813773
$ws.Range("J19").Value = 1.031198718048999
$ws.Range("K19").Value = 1.031575466292568
$ws.Range("L19").Value = 1.041425134870662
$ws.Range("M19").Value = 1.047758424175623
$ws.Range("N19").Value = 1.014385853890475

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024717622788564
$ws.Range("D20").Value = 1.028007378042642
$ws.Range("E20").Value = 1.037784010339346
$ws.Range("F20").Value = 1.044089354253717
$ws.Range("I20").Value = 1.029144560420518
$ws.Range("J20").Value = 1.030990670656119
$ws.Range("K20").Value = 1.031409080033153
$ws.Range("L20").Value = 1.041151330212333
$ws.Range("M20").Value = 1.047434875584943
$ws.Range("N20").Value = 1.014315822388818

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02363796576773
$ws.Range("D21").Value = 1.027254301384515
$ws.Range("E21").Value = 1.036683874404787
$ws.Range("F21").Value = 1.042828842929846
$ws.Range("I21").Value = 1.029008053328472
$ws.Range("J21").Value = 1.03031403310332
$ws.Range("K21").Value = 1.030867294517998
$ws.Range("L21").Value = 1.040261469278354
$ws.Range("M21").Value = 1.046383750832176
$ws.Range("N21").Value = 1.014087956806498

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022959578618477
$ws.Range("D22").Value = 1.026780902907285
$ws.Range("E22").Value = 1.035993069832127
$ws.Range("F22").Value = 1.042037520093114
$ws.Range("I22").Value = 1.0289207863034
$ws.Range("J22").Value = 1.029888411982076
$ws.Range("K22").Value = 1.03052600698418
$ws.Range("L22").Value = 1.039702215470702
$ws.Range("M22").Value = 1.045723455228378
$ws.Range("N22").Value = 1.013944547307774

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023319162957098
$ws.Range("D23").Value = 1.027031851380886
$ws.Range("E23").Value = 1.036359193334168
$ws.Range("F23").Value = 1.042456900331826
$ws.Range("I23").Value = 1.028967184511526
$ws.Range("J23").Value = 1.030114059786251
$ws.Range("K23").Value = 1.030706990869219
$ws.Range("L23").Value = 1.039998663783694
$ws.Range("M23").Value = 1.046073434441112
$ws.Range("N23").Value = 1.014020584698052

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024735730732877
$ws.Range("D24").Value = 1.028020004917994
$ws.Range("E24").Value = 1.037802469367765
$ws.Range("F24").Value = 1.044110507370482
$ws.Range("I24").Value = 1.029146824542651
$ws.Range("J24").Value = 1.031002011265847
$ws.Range("K24").Value = 1.031418152135756
$ws.Range("L24").Value = 1.041166252815574
$ws.Range("M24").Value = 1.047452507754781
$ws.Range("N24").Value = 1.014319640167559

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02638145851282
$ws.Range("D25").Value = 1.029167054890986
$ws.Range("E25").Value = 1.03948117800152
$ws.Range("F25").Value = 1.046034676599259
$ws.Range("I25").Value = 1.029349003185317
$ws.Range("J25").Value = 1.032031569631497
$ws.Range("K25").Value = 1.032240576589267
$ws.Range("L25").Value = 1.042522183152207
$ws.Range("M25").Value = 1.04905538480927
$ws.Range("N25").Value = 1.014666051977376
